# Add "2022-Q1" fund-holdings sheet before the "总计" (totals) sheet,
# and prepend a matching summary row to the "总计" (totals) sheet.

$wb = $excel.ActiveWorkbook

# --- 1. Create the new "2022-Q1" sheet by duplicating the "2021-Q4" sheet
#        (same layout/styles) and placing the copy right before "总计". ---
$template = $wb.Worksheets.Item("2021-Q4")
$totalSheet = $wb.Worksheets.Item("总计")
$template.Copy($totalSheet)

# The copy is inserted immediately before "总计" and auto-named "2021-Q4 (2)".
$newSheet = $wb.Worksheets.Item("2021-Q4 (2)")
$newSheet.Name = "2022-Q1"

# --- 2. Overwrite the copied data with the real 2022-Q1 fund holdings. ---
# Columns B-G hold text values (fund code / name / amounts as strings);
# force text storage so numeric-looking strings aren't reinterpreted.
$newSheet.Range("B2:G3").NumberFormat = "@"

$newSheet.Range("B2").Value = "010387"
$newSheet.Range("C2").Value = "易方达医药生物股票A"
$newSheet.Range("D2").Value = "22.51"
$newSheet.Range("E2").Value = "83.49"
$newSheet.Range("F2").Value = "3.91"
$newSheet.Range("G2").Value = "0.8801"
$newSheet.Range("H2").Value = 6

$newSheet.Range("B3").Value = "010388"
$newSheet.Range("C3").Value = "易方达医药生物股票C"
$newSheet.Range("D3").Value = "4.78"
$newSheet.Range("E3").Value = "83.49"
$newSheet.Range("F3").Value = "3.91"
$newSheet.Range("G3").Value = "0.1869"
$newSheet.Range("H3").Value = 6

# --- 3. Insert a new row at the top of the "总计" sheet's data (row 2),
#        pushing the existing quarterly summaries down one row. ---
# (Re-fetch the "总计" sheet: references captured before the sheet copy
#  above are no longer reliable once the sheet collection changes.)
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()

# The inserted row picks up stray formatting from the header row above it;
# strip that so the new row starts from a clean/default style.
$totalSheet.Range("A2:D2").ClearFormats()

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 1.07

# Restore the index-column style (A2 should match A3:A6's "s=2" look).
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

# --- 4. Renumber the index column (A) for the rows pushed down. ---
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4
